$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows above row 3 (pushes the old row3..row11 block down to row6..row14),
# inheriting the formatting (wrap text + left align) of row 3's cells.
$ws.Rows("3:5").Insert()

# --- Fill in the new trading-plan entry for 20170103 (Tuesday) in the freshly inserted row 3 ---
$ws.Range("A3").Value = 20170103
$ws.Range("B3").Value = "Tuesday"
$ws.Range("C3").Value = 20170103
$ws.Range("D3").Value = "Tuesday"
$ws.Range("E3").Value = "huge gap down following the new year's long weekend, due to warmer than normal weather forecast in mid Jan. The Ng price may swing up a little more likely for a bull trap in the early morning (like now, around 1030am) and then will resume down trend and continue this to tomorrow, the day before the report day, which is expected to be ~-70, quite bearish but it's probabl6y already priced in by now. The next daily susport area is 3.315 -3.341"
$ws.Range("F3").Value = "be patient to see -10% in NG today! "
$ws.Range("G3").Value = "20170103-15min and 20170103-1h"

# C3 picked up the same (wrap-only, no horizontal align) style as the date cells in column A.
$ws.Range("C3").Style = $ws.Range("A3").Style

# Row 3 grows taller to fit the wrapped note text.
$ws.Rows("3").RowHeight = 72

# Column G now has the standard default column width/style like the other columns.
$ws.Columns("G").ColumnWidth = 8.88671875

# Scroll the view over so column C is the first visible column (as in the saved workbook).
$ws.Application.ActiveWindow.ScrollColumn = 3
